# Calendar time slots - first integration with guest house, basic function working
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Time Slot"
$ws.Range("C1").Value = "Display Time"
$ws.Range("D1").Value = "Capacity"
$ws.Range("E1").Value = "Booked"
$ws.Range("F1").Value = "Available"

# Data rows - 2024-03-20, slots 15:00 / 16:00 / 17:00
$dates = @("2024-03-20", "2024-03-20", "2024-03-20")
$times = @(0.625, 0.66666666666666663, 0.70833333333333337)
$capacity = 10
$booked = 0

for ($i = 0; $i -lt 3; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $times[$i]
    $ws.Cells.Item($row, 3).Value = $times[$i]
    $ws.Cells.Item($row, 4).Value = $capacity
    $ws.Cells.Item($row, 5).Value = $booked
    $ws.Cells.Item($row, 6).Value = $capacity - $booked
}

# Number formats
$ws.Range("A1:F1").HorizontalAlignment = -4131
$ws.Range("A2:A4").NumberFormat = "mm-dd-yy"
$ws.Range("B2:B4").NumberFormat = "h:mm"
$ws.Range("C2:C4").NumberFormat = "h:mm AM/PM"

$ws.Range("A1:F4").HorizontalAlignment = -4131

# Column widths - best fit
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null

$ws.PageSetup.OddFooter = "&C_x000D_&1#&`"Calibri`"&8&K000000 Internal"

$ws.Range("H10").Select()
